# [Fonds de solidarite] Add 2020-09-07 data
# Update nombre_aides (col C) and montant_total (col D) for the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 3).Value = "'1052"
$ws.Cells.Item(3, 4).Value = "'3439244.33"
$ws.Cells.Item(4, 3).Value = "'430"
$ws.Cells.Item(4, 4).Value = "'1815645.25"
$ws.Cells.Item(5, 3).Value = "'120"
$ws.Cells.Item(5, 4).Value = "'576128.09"
$ws.Cells.Item(6, 3).Value = "'32"
$ws.Cells.Item(6, 4).Value = "'208643.82"
$ws.Cells.Item(22, 3).Value = "'334"
$ws.Cells.Item(22, 4).Value = "'1019762.20"
$ws.Cells.Item(23, 3).Value = "'122"
$ws.Cells.Item(23, 4).Value = "'513160.00"
$ws.Cells.Item(46, 3).Value = "'98"
$ws.Cells.Item(46, 4).Value = "'441274.61"
$ws.Cells.Item(47, 3).Value = "'54"
$ws.Cells.Item(47, 4).Value = "'311703.00"
$ws.Cells.Item(48, 3).Value = "'31"
$ws.Cells.Item(48, 4).Value = "'213697.00"
$ws.Cells.Item(49, 3).Value = "'6"
$ws.Cells.Item(49, 4).Value = "'36000.00"
$ws.Cells.Item(50, 3).Value = "'17"
$ws.Cells.Item(50, 4).Value = "'37850.00"
$ws.Cells.Item(52, 3).Value = "'622"
$ws.Cells.Item(52, 4).Value = "'2282958.21"
$ws.Cells.Item(53, 3).Value = "'273"
$ws.Cells.Item(53, 4).Value = "'1252878.76"
$ws.Cells.Item(55, 3).Value = "'28"
$ws.Cells.Item(55, 4).Value = "'163213.00"
$ws.Cells.Item(56, 3).Value = "'27"
$ws.Cells.Item(56, 4).Value = "'84220.65"
$ws.Cells.Item(60, 3).Value = "'635"
$ws.Cells.Item(60, 4).Value = "'3125291.45"
$ws.Cells.Item(84, 3).Value = "'911"
$ws.Cells.Item(84, 4).Value = "'2945833.26"
$ws.Cells.Item(101, 3).Value = "'311"
$ws.Cells.Item(101, 4).Value = "'821687.37"
$ws.Cells.Item(102, 3).Value = "'1280"
$ws.Cells.Item(102, 4).Value = "'4045488.43"
$ws.Cells.Item(103, 3).Value = "'478"
$ws.Cells.Item(103, 4).Value = "'2018990.62"
$ws.Cells.Item(104, 3).Value = "'127"
$ws.Cells.Item(104, 4).Value = "'611996.00"
$ws.Cells.Item(105, 3).Value = "'38"
$ws.Cells.Item(105, 4).Value = "'248657.00"
$ws.Cells.Item(106, 3).Value = "'74"
$ws.Cells.Item(106, 4).Value = "'159956.16"
